$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# --- Add 4 new rows (18-21) for the menu pushbuttons ---
# (do this before touching row 17's own formatting, since rows 19/21 need
# row 17's *original* B/D style, i.e. s="8")

# Row 18: style like row 16 (A/C = s5, B/D = s6)
$ws.Range("A16:D16").Copy() | Out-Null
$ws.Range("A18:D18").PasteSpecial($xlPasteFormats) | Out-Null

# Row 19: style like row 17 (A/C = s3, B/D = s8)
$ws.Range("A17:D17").Copy() | Out-Null
$ws.Range("A19:D19").PasteSpecial($xlPasteFormats) | Out-Null

# Row 20: style like row 16 (A/C = s5, B/D = s6)
$ws.Range("A16:D16").Copy() | Out-Null
$ws.Range("A20:D20").PasteSpecial($xlPasteFormats) | Out-Null

# Row 21: style like row 17 (A/C = s3, B/D = s8)
$ws.Range("A17:D17").Copy() | Out-Null
$ws.Range("A21:D21").PasteSpecial($xlPasteFormats) | Out-Null

# --- Fix B17/D17 cell style (was s="8", should become s="4") ---
# B5/D5 already use style s="4" so copy their formatting onto B17/D17.
$ws.Range("B5").Copy() | Out-Null
$ws.Range("B17").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("D5").Copy() | Out-Null
$ws.Range("D17").PasteSpecial($xlPasteFormats) | Out-Null

$excel.CutCopyMode = 0

# --- Fill in values for the new rows ---
# Shared-string table order follows column-by-column entry (all of column B,
# then all of column D), so write the cells in that same order.
$ws.Range("A18").Value = 15
$ws.Range("C18").Value = 8
$ws.Range("A19").Value = 16
$ws.Range("C19").Value = 9
$ws.Range("A20").Value = 17
$ws.Range("C20").Value = 10
$ws.Range("A21").Value = 18
$ws.Range("C21").Value = 11

$ws.Range("B18").Value = "Pushbutton vermelho"
$ws.Range("B19").Value = "Pushbutton verde"
$ws.Range("B20").Value = "Pushbutton amarelo"
$ws.Range("B21").Value = "Pushbutton azul"

$ws.Range("D18").Value = "Botão para seleção no menu de opções"
$ws.Range("D19").Value = "Botão para retornar a tela inicial do relógio"
$ws.Range("D20").Value = "Botão para seleção de subitem do menu de opções"
$ws.Range("D21").Value = "Botão para entrada de dados de configuração"

# --- Update selection to reflect where the user ended up after editing ---
$ws.Range("H22").Select() | Out-Null
